$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new columns before column D (quarters shift right by two columns).
$ws.Range("D1:E1").EntireColumn.Insert()

# New D:E columns should inherit the formatting that column F (the former
# column D) carries, so copy formats across before writing values.
$ws.Range("F1:F102").Copy()
$ws.Range("D1:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Per-row data for the two newly inserted quarters (columns D and E).
$newQuarterData = @(
    @{Row=7; D=43465; E=43372}
    @{Row=8; D=1662000; E=1651000}
    @{Row=9; D=1549000; E=1529000}
    @{Row=10; D=113000; E=122000}
    @{Row=11; D=$null; E=$null}
    @{Row=12; D="NA"; E="NA"}
    @{Row=13; D=0; E=0}
    @{Row=14; D=0; E=0}
    @{Row=15; D=0; E=0}
    @{Row=16; D=$null; E=$null}
    @{Row=17; D=1619000; E=1614000}
    @{Row=18; D=43000; E=37000}
    @{Row=19; D=$null; E=$null}
    @{Row=20; D=-171000; E=35000}
    @{Row=21; D=-95000; E=106000}
    @{Row=22; D=13000; E=12000}
    @{Row=23; D=-141000; E=60000}
    @{Row=24; D=-65000; E=26000}
    @{Row=25; D=0; E=0}
    @{Row=26; D=-76000; E=34000}
    @{Row=27; D=-77000; E=35000}
    @{Row=28; D=0; E=0}
    @{Row=29; D=-14000; E=0}
    @{Row=30; D=0; E=0}
    @{Row=31; D=0; E=0}
    @{Row=32; D=171000; E=-35000}
    @{Row=33; D=-91000; E=35000}
    @{Row=34; D=0; E=0}
    @{Row=35; D=-91000; E=35000}
    @{Row=38; D=43465; E=43372}
    @{Row=39; D=$null; E=$null}
    @{Row=40; D=$null; E=$null}
    @{Row=41; D=194000; E=78000}
    @{Row=42; D=1336000; E=1590000}
    @{Row=43; D=551000; E=506000}
    @{Row=44; D=815000; E=874000}
    @{Row=45; D=131000; E=124000}
    @{Row=46; D=3027000; E=3172000}
    @{Row=47; D=804000; E=823000}
    @{Row=48; D=1160000; E=1098000}
    @{Row=49; D=236000; E="NA"}
    @{Row=50; D=0; E=0}
    @{Row=51; D=0; E=0}
    @{Row=52; D=80000; E=342000}
    @{Row=53; D=0; E=0}
    @{Row=54; D=5307000; E=5435000}
    @{Row=55; D=$null; E=$null}
    @{Row=56; D=$null; E=$null}
    @{Row=57; D=238000; E=227000}
    @{Row=58; D=187000; E=172000}
    @{Row=59; D=359000; E=370000}
    @{Row=60; D=784000; E=769000}
    @{Row=61; D=739000; E=741000}
    @{Row=62; D=455000; E=499000}
    @{Row=63; D=0; E=0}
    @{Row=64; D=0; E=0}
    @{Row=65; D=0; E=0}
    @{Row=66; D=1989000; E=2023000}
    @{Row=67; D=$null; E=$null}
    @{Row=68; D=0; E=0}
    @{Row=69; D=0; E=0}
    @{Row=70; D=0; E=0}
    @{Row=71; D=0; E=0}
    @{Row=72; D=3727000; E=3825000}
    @{Row=73; D=0; E=0}
    @{Row=74; D=0; E=0}
    @{Row=75; D=0; E=0}
    @{Row=76; D=3318000; E=3412000}
    @{Row=77; D=0; E=0}
    @{Row=80; D=43465; E=43372}
    @{Row=81; D=-91000; E=35000}
    @{Row=82; D=$null; E=$null}
    @{Row=83; D=33000; E=34000}
    @{Row=84; D=0; E=0}
    @{Row=85; D=0; E=0}
    @{Row=86; D=0; E=0}
    @{Row=87; D=0; E=0}
    @{Row=88; D=0; E=0}
    @{Row=89; D=72000; E=102000}
    @{Row=90; D=$null; E=$null}
    @{Row=91; D=-67000; E=-37000}
    @{Row=92; D=0; E=0}
    @{Row=93; D=0; E=0}
    @{Row=94; D=37000; E=-340000}
    @{Row=95; D=$null; E=$null}
    @{Row=96; D=-1000; E=-2000}
    @{Row=97; D=0; E=0}
    @{Row=98; D=0; E=0}
    @{Row=99; D=0; E=0}
    @{Row=100; D=5000; E=240000}
    @{Row=101; D=2000; E=-3000}
    @{Row=102; D=116000; E=-1000}
)

foreach ($item in $newQuarterData) {
    $r = $item.Row
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
}
